# Apply the "Sasank" merge edit to the User Stories sheet.
# - Replace the "track my food expense" user story (row 2) with the
#   "upload an image of food" / Recipe Recommendation story.
# - Replace the "easily upload my grocery receipts" user story (row 3) with
#   the "see detailed nutritional information" / Nutritional Analysis story.
# - Update the sheet view zoom / selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Stories")

# Values are assigned in this precise order so that new shared-string table
# entries are appended in the same sequence the original authoring session
# produced them in.
$ws.Range("C2").Value = "Student "
$ws.Range("D2").Value = "upload an image of food"
$ws.Range("C3").Value = "Parent"
$ws.Range("D3").Value = "see detailed nutritional information for the food in the uploaded image"
$ws.Range("E3").Value = "I can make informed dietary choices"
$ws.Range("E2").Value = "I can receive video recipes related to the image"
$ws.Range("F3").Value = "Nutritional Analysis"
$ws.Range("F2").Value = "Recipe Recommendation"

# --- Sheet view: zoom + active selection ----------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 125
$ws.Range("D21").Select()
